$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.130.45'
$ws.Range("E2").Value = '  +1.21%  '
$ws.Range("D3").Value = '2.382.84'
$ws.Range("E3").Value = '  +3.62%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = "'302.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.41%  '
$ws.Range("D6").Value = "'97.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.46%  '
$ws.Range("D7").Value = "'0.504"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  +1.48%  '
$ws.Range("D10").Value = "'34.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.69%  '
$ws.Range("D11").Value = "'0.0788"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.60%  '
$ws.Range("D12").Value = "'0.122"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.34%  '
$ws.Range("D13").Value = "'18.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.75%  '
$ws.Range("D14").Value = "'6.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.71%  '
$ws.Range("D15").Value = '2.753.43'
$ws.Range("E15").Value = '  +3.57%  '
$ws.Range("D16").Value = '2.346.08'
$ws.Range("E16").Value = '  +2.21%  '
$ws.Range("D17").Value = "'0.807"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.59%  '
$ws.Range("D18").Value = '43.139.52'
$ws.Range("E18").Value = '  +1.40%  '
$ws.Range("E19").Value = '  -1.91%  '
$ws.Range("D20").Value = "'6.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.71%  '
$ws.Range("D21").Value = '0.0₃0887'
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("D22").Value = "'68.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.04%  '
$ws.Range("D23").Value = "'235.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").Value = "'2.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.46%  '
$ws.Range("E25").Value = '  +1.70%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").Value = "'24.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.77%  '
$ws.Range("E28").Value = '  +0.39%  '
$ws.Range("D29").Value = "'9.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.88%  '
$ws.Range("D30").Value = "'31.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.07%  '
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").Value = "'5.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.17%  '
$ws.Range("D33").Value = "'0.0750"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.07%  '
$ws.Range("E34").Value = '  -1.03%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = "'1.86"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.15%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = "'0.105"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.22%  '
$ws.Range("D37").Value = "'2.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.33%  '
$ws.Range("D38").Value = "'4.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.35%  '
$ws.Range("D39").Value = "'2.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.81%  '
$ws.Range("D40").Value = "'22.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +11.61%  '
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("D42").Value = "'106.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -35.28%  '
$ws.Range("D43").Value = '1.959.59'
$ws.Range("E43").Value = '  +0.55%  '
$ws.Range("E44").Value = '  +0.79%  '
$ws.Range("E45").Value = '  +2.10%  '
$ws.Range("E46").Value = '  +0.57%  '
$ws.Range("D47").Value = "'9.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -11.68%  '
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").Value = '2.604.06'
$ws.Range("E48").Value = '  +3.04%  '
$ws.Range("B49").Value = 'MultiversX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D49").Value = "'52.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.56%  '
$ws.Range("D50").Value = "'1.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.94%  '
$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").Value = "'71.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.24%  '
